$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1580.7142  # ALC!H19 (1584.2858 -> 1580.7142)
$ws.Cells.Item(19, 9).Value = 1414.75  # ALC!I19 (1427.25 -> 1414.75)
$ws.Cells.Item(19, 11).Value = 1414.75  # ALC!K19 (1427.25 -> 1414.75)
$ws.Cells.Item(19, 13).Value = -1239.75  # ALC!M19 (-1252.25 -> -1239.75)
$ws.Cells.Item(33, 8).Value = 4137.125  # ALC!H33 (4144.6 -> 4137.125)
$ws.Cells.Item(33, 9).Value = 2623.4  # ALC!I33 (2252.6667 -> 2623.4)
$ws.Cells.Item(33, 10).Value = 6660  # ALC!J33 (6982.5 -> 6660)
$ws.Cells.Item(33, 11).Value = 2623.4  # ALC!K33 (2252.6667 -> 2623.4)
$ws.Cells.Item(33, 12).Value = 6660  # ALC!L33 (6982.5 -> 6660)
$ws.Cells.Item(33, 13).Value = -2394.4  # ALC!M33 (-2023.6667 -> -2394.4)
$ws.Cells.Item(33, 14).Value = -7118  # ALC!N33 (-7440.5 -> -7118)
$ws.Cells.Item(55, 8).Value = 395  # ALC!H55 (450 -> 395)
$ws.Cells.Item(55, 9).Value = 395  # ALC!I55 (0 -> 395)
$ws.Cells.Item(55, 10).Value = 0  # ALC!J55 (450 -> 0)
$ws.Cells.Item(55, 11).Value = 395  # ALC!K55 (0 -> 395)
$ws.Cells.Item(55, 12).Value = 0  # ALC!L55 (450 -> 0)
$ws.Cells.Item(55, 13).Value = -181  # ALC!M55 (None -> -181)
$ws.Cells.Item(55, 14).Value = ""  # ALC!N55 clear (was -878)
$ws.Cells.Item(76, 8).Value = 7366.6665  # ALC!H76 (5696.6 -> 7366.6665)
$ws.Cells.Item(76, 9).Value = 3200  # ALC!I76 (3150 -> 3200)
$ws.Cells.Item(76, 10).Value = 8200  # ALC!J76 (6333.25 -> 8200)
$ws.Cells.Item(76, 11).Value = 3200  # ALC!K76 (3150 -> 3200)
$ws.Cells.Item(76, 12).Value = 8200  # ALC!L76 (6333.25 -> 8200)
$ws.Cells.Item(76, 13).Value = -2885  # ALC!M76 (-2835 -> -2885)
$ws.Cells.Item(76, 14).Value = -8830  # ALC!N76 (-6963.25 -> -8830)
$ws.Cells.Item(79, 8).Value = 7366.6665  # ALC!H79 (5696.6 -> 7366.6665)
$ws.Cells.Item(79, 9).Value = 3200  # ALC!I79 (3150 -> 3200)
$ws.Cells.Item(79, 10).Value = 8200  # ALC!J79 (6333.25 -> 8200)
$ws.Cells.Item(79, 11).Value = 3200  # ALC!K79 (3150 -> 3200)
$ws.Cells.Item(79, 12).Value = 8200  # ALC!L79 (6333.25 -> 8200)
$ws.Cells.Item(79, 13).Value = -2108  # ALC!M79 (-2058 -> -2108)
$ws.Cells.Item(79, 14).Value = -10384  # ALC!N79 (-8517.25 -> -10384)
$ws.Cells.Item(116, 8).Value = 3567.7778  # ALC!H116 (3451.625 -> 3567.7778)
$ws.Cells.Item(116, 10).Value = 4832.2  # ALC!J116 (4916 -> 4832.2)
$ws.Cells.Item(116, 12).Value = 4832.2  # ALC!L116 (4916 -> 4832.2)
$ws.Cells.Item(116, 14).Value = -11716.2  # ALC!N116 (-11800 -> -11716.2)

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1378.8  # ARM!H5 (1315.3334 -> 1378.8)
$ws.Cells.Item(5, 9).Value = 1223.5  # ARM!I5 (1178.4 -> 1223.5)
$ws.Cells.Item(5, 11).Value = 1223.5  # ARM!K5 (1178.4 -> 1223.5)
$ws.Cells.Item(5, 13).Value = -1111.5  # ARM!M5 (-1066.4 -> -1111.5)
$ws.Cells.Item(61, 8).Value = 5244.6  # ARM!H61 (7266.6665 -> 5244.6)
$ws.Cells.Item(61, 9).Value = 6002.75  # ARM!I61 (7266.6665 -> 6002.75)
$ws.Cells.Item(61, 10).Value = 2212  # ARM!J61 (0 -> 2212)
$ws.Cells.Item(61, 11).Value = 6002.75  # ARM!K61 (7266.6665 -> 6002.75)
$ws.Cells.Item(61, 12).Value = 2212  # ARM!L61 (0 -> 2212)
$ws.Cells.Item(61, 13).Value = -5790.75  # ARM!M61 (-7054.6665 -> -5790.75)
$ws.Cells.Item(61, 14).Value = -2636  # ARM!N61 (None -> -2636)
$ws.Cells.Item(110, 8).Value = 1367.75  # ARM!H110 (1317.1538 -> 1367.75)
$ws.Cells.Item(110, 9).Value = 1367.75  # ARM!I110 (1389.4166 -> 1367.75)
$ws.Cells.Item(110, 10).Value = 0  # ARM!J110 (450 -> 0)
$ws.Cells.Item(110, 11).Value = 1367.75  # ARM!K110 (1389.4166 -> 1367.75)
$ws.Cells.Item(110, 12).Value = 0  # ARM!L110 (450 -> 0)
$ws.Cells.Item(110, 13).Value = 677.25  # ARM!M110 (655.5834 -> 677.25)
$ws.Cells.Item(110, 14).Value = ""  # ARM!N110 clear (was -4540)
$ws.Cells.Item(122, 8).Value = 5268.8  # ARM!H122 (4797.4443 -> 5268.8)
$ws.Cells.Item(122, 9).Value = 5026.091  # ARM!I122 (4472.0713 -> 5026.091)
$ws.Cells.Item(122, 11).Value = 15078.273  # ARM!K122 (13416.2139 -> 15078.273)
$ws.Cells.Item(122, 13).Value = -12628.273  # ARM!M122 (-10966.2139 -> -12628.273)
$ws.Cells.Item(136, 8).Value = 5244.6  # ARM!H136 (7266.6665 -> 5244.6)
$ws.Cells.Item(136, 9).Value = 6002.75  # ARM!I136 (7266.6665 -> 6002.75)
$ws.Cells.Item(136, 10).Value = 2212  # ARM!J136 (0 -> 2212)
$ws.Cells.Item(136, 11).Value = 18008.25  # ARM!K136 (21799.9995 -> 18008.25)
$ws.Cells.Item(136, 12).Value = 6636  # ARM!L136 (0 -> 6636)
$ws.Cells.Item(136, 13).Value = -15458.25  # ARM!M136 (-19249.9995 -> -15458.25)
$ws.Cells.Item(136, 14).Value = -11736  # ARM!N136 (None -> -11736)

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1378.8  # BSM!H4 (1315.3334 -> 1378.8)
$ws.Cells.Item(4, 9).Value = 1223.5  # BSM!I4 (1178.4 -> 1223.5)
$ws.Cells.Item(4, 11).Value = 1223.5  # BSM!K4 (1178.4 -> 1223.5)
$ws.Cells.Item(4, 13).Value = -1108.5  # BSM!M4 (-1063.4 -> -1108.5)
$ws.Cells.Item(22, 8).Value = 986.9231  # BSM!H22 (987.7692 -> 986.9231)
$ws.Cells.Item(22, 9).Value = 927.8570999999999  # BSM!I22 (929.4286 -> 927.8570999999999)
$ws.Cells.Item(22, 11).Value = 927.8570999999999  # BSM!K22 (929.4286 -> 927.8570999999999)
$ws.Cells.Item(22, 13).Value = -754.8570999999999  # BSM!M22 (-756.4286 -> -754.8570999999999)
$ws.Cells.Item(35, 8).Value = 92000  # BSM!H35 (72000 -> 92000)
$ws.Cells.Item(35, 10).Value = 92000  # BSM!J35 (72000 -> 92000)
$ws.Cells.Item(35, 12).Value = 92000  # BSM!L35 (72000 -> 92000)
$ws.Cells.Item(35, 14).Value = -92620  # BSM!N35 (-72620 -> -92620)
$ws.Cells.Item(94, 8).Value = 670.59375  # BSM!H94 (688.3 -> 670.59375)
$ws.Cells.Item(94, 9).Value = 667.7  # BSM!I94 (681.7895 -> 667.7)
$ws.Cells.Item(94, 10).Value = 675.4167  # BSM!J94 (699.5454999999999 -> 675.4167)
$ws.Cells.Item(94, 11).Value = 667.7  # BSM!K94 (681.7895 -> 667.7)
$ws.Cells.Item(94, 12).Value = 675.4167  # BSM!L94 (699.5454999999999 -> 675.4167)
$ws.Cells.Item(94, 13).Value = -216.7  # BSM!M94 (-230.7895 -> -216.7)
$ws.Cells.Item(94, 14).Value = -1577.4167  # BSM!N94 (-1601.5455 -> -1577.4167)

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2766.8928  # CRP!H58 (2724.8965 -> 2766.8928)
$ws.Cells.Item(58, 9).Value = 2808  # CRP!I58 (2694.8333 -> 2808)
$ws.Cells.Item(58, 10).Value = 2740.2942  # CRP!J58 (2746.1177 -> 2740.2942)
$ws.Cells.Item(58, 11).Value = 2808  # CRP!K58 (2694.8333 -> 2808)
$ws.Cells.Item(58, 12).Value = 2740.2942  # CRP!L58 (2746.1177 -> 2740.2942)
$ws.Cells.Item(58, 13).Value = -2605  # CRP!M58 (-2491.8333 -> -2605)
$ws.Cells.Item(58, 14).Value = -3146.2942  # CRP!N58 (-3152.1177 -> -3146.2942)
$ws.Cells.Item(107, 8).Value = 1597.2285  # CRP!H107 (1641.6177 -> 1597.2285)
$ws.Cells.Item(107, 9).Value = 1430.1786  # CRP!I107 (1446.6428 -> 1430.1786)
$ws.Cells.Item(107, 10).Value = 2265.4285  # CRP!J107 (2551.5 -> 2265.4285)
$ws.Cells.Item(107, 11).Value = 1430.1786  # CRP!K107 (1446.6428 -> 1430.1786)
$ws.Cells.Item(107, 12).Value = 2265.4285  # CRP!L107 (2551.5 -> 2265.4285)
$ws.Cells.Item(107, 13).Value = 489.8214  # CRP!M107 (473.3571999999999 -> 489.8214)
$ws.Cells.Item(107, 14).Value = -6105.4285  # CRP!N107 (-6391.5 -> -6105.4285)
$ws.Cells.Item(132, 8).Value = 7478.25  # CRP!H132 (5508 -> 7478.25)
$ws.Cells.Item(132, 9).Value = 0  # CRP!I132 (1433.4 -> 0)
$ws.Cells.Item(132, 10).Value = 7478.25  # CRP!J132 (9582.6 -> 7478.25)
$ws.Cells.Item(132, 11).Value = 0  # CRP!K132 (4300.200000000001 -> 0)
$ws.Cells.Item(132, 12).Value = 22434.75  # CRP!L132 (28747.8 -> 22434.75)
$ws.Cells.Item(132, 13).Value = ""  # CRP!M132 clear (was -1770.200000000001)
$ws.Cells.Item(132, 14).Value = -27494.75  # CRP!N132 (-33807.8 -> -27494.75)
$ws.Cells.Item(134, 8).Value = 6000  # CRP!H134 (6345.75 -> 6000)
$ws.Cells.Item(134, 9).Value = 5038.077  # CRP!I134 (5333.522 -> 5038.077)
$ws.Cells.Item(134, 11).Value = 15114.231  # CRP!K134 (16000.566 -> 15114.231)
$ws.Cells.Item(134, 13).Value = -12579.231  # CRP!M134 (-13465.566 -> -12579.231)
$ws.Cells.Item(136, 8).Value = 2766.8928  # CRP!H136 (2724.8965 -> 2766.8928)
$ws.Cells.Item(136, 9).Value = 2808  # CRP!I136 (2694.8333 -> 2808)
$ws.Cells.Item(136, 10).Value = 2740.2942  # CRP!J136 (2746.1177 -> 2740.2942)
$ws.Cells.Item(136, 11).Value = 8424  # CRP!K136 (8084.499899999999 -> 8424)
$ws.Cells.Item(136, 12).Value = 8220.882599999999  # CRP!L136 (8238.3531 -> 8220.882599999999)
$ws.Cells.Item(136, 13).Value = -5874  # CRP!M136 (-5534.499899999999 -> -5874)
$ws.Cells.Item(136, 14).Value = -13320.8826  # CRP!N136 (-13338.3531 -> -13320.8826)

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(25, 8).Value = 1777.5  # CUL!H25 (0 -> 1777.5)
$ws.Cells.Item(25, 9).Value = 555  # CUL!I25 (0 -> 555)
$ws.Cells.Item(25, 10).Value = 3000  # CUL!J25 (0 -> 3000)
$ws.Cells.Item(25, 11).Value = 1665  # CUL!K25 (0 -> 1665)
$ws.Cells.Item(25, 12).Value = 9000  # CUL!L25 (0 -> 9000)
$ws.Cells.Item(25, 13).Value = -1496  # CUL!M25 (None -> -1496)
$ws.Cells.Item(25, 14).Value = -9338  # CUL!N25 (None -> -9338)
$ws.Cells.Item(30, 8).Value = 1777.5  # CUL!H30 (0 -> 1777.5)
$ws.Cells.Item(30, 9).Value = 555  # CUL!I30 (0 -> 555)
$ws.Cells.Item(30, 10).Value = 3000  # CUL!J30 (0 -> 3000)
$ws.Cells.Item(30, 11).Value = 1665  # CUL!K30 (0 -> 1665)
$ws.Cells.Item(30, 12).Value = 9000  # CUL!L30 (0 -> 9000)
$ws.Cells.Item(30, 13).Value = -1563  # CUL!M30 (None -> -1563)
$ws.Cells.Item(30, 14).Value = -9204  # CUL!N30 (None -> -9204)
$ws.Cells.Item(86, 8).Value = 665.1818  # CUL!H86 (634.75 -> 665.1818)
$ws.Cells.Item(86, 10).Value = 1118  # CUL!J86 (981.6667 -> 1118)
$ws.Cells.Item(86, 12).Value = 3354  # CUL!L86 (2945.0001 -> 3354)
$ws.Cells.Item(86, 14).Value = -5726  # CUL!N86 (-5317.0001 -> -5726)
$ws.Cells.Item(89, 8).Value = 665.1818  # CUL!H89 (634.75 -> 665.1818)
$ws.Cells.Item(89, 10).Value = 1118  # CUL!J89 (981.6667 -> 1118)
$ws.Cells.Item(89, 12).Value = 10062  # CUL!L89 (8835.0003 -> 10062)
$ws.Cells.Item(89, 14).Value = -21918  # CUL!N89 (-20691.0003 -> -21918)
$ws.Cells.Item(131, 8).Value = 3045.4827  # CUL!H131 (3052.3794 -> 3045.4827)
$ws.Cells.Item(131, 10).Value = 3136.5  # CUL!J131 (3145.5908 -> 3136.5)
$ws.Cells.Item(131, 12).Value = 9409.5  # CUL!L131 (9436.7724 -> 9409.5)
$ws.Cells.Item(131, 14).Value = -19489.5  # CUL!N131 (-19516.7724 -> -19489.5)
$ws.Cells.Item(134, 8).Value = 1239.1666  # CUL!H134 (1732.5 -> 1239.1666)
$ws.Cells.Item(134, 9).Value = 1239.1666  # CUL!I134 (1732.5 -> 1239.1666)
$ws.Cells.Item(134, 11).Value = 3717.4998  # CUL!K134 (5197.5 -> 3717.4998)
$ws.Cells.Item(134, 13).Value = 1352.5002  # CUL!M134 (-127.5 -> 1352.5002)
$ws.Cells.Item(139, 8).Value = 5789.6  # CUL!H139 (4063.6 -> 5789.6)
$ws.Cells.Item(139, 9).Value = 5789.6  # CUL!I139 (4063.6 -> 5789.6)
$ws.Cells.Item(139, 11).Value = 17368.8  # CUL!K139 (12190.8 -> 17368.8)
$ws.Cells.Item(139, 13).Value = -12228.8  # CUL!M139 (-7050.799999999999 -> -12228.8)

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2608.0435  # GSM!H80 (2662.9546 -> 2608.0435)
$ws.Cells.Item(80, 9).Value = 2439.0667  # GSM!I80 (2513.2856 -> 2439.0667)
$ws.Cells.Item(80, 11).Value = 2439.0667  # GSM!K80 (2513.2856 -> 2439.0667)
$ws.Cells.Item(80, 13).Value = -1441.0667  # GSM!M80 (-1515.2856 -> -1441.0667)
$ws.Cells.Item(83, 8).Value = 2608.0435  # GSM!H83 (2662.9546 -> 2608.0435)
$ws.Cells.Item(83, 9).Value = 2439.0667  # GSM!I83 (2513.2856 -> 2439.0667)
$ws.Cells.Item(83, 11).Value = 12195.3335  # GSM!K83 (12566.428 -> 12195.3335)
$ws.Cells.Item(83, 13).Value = -7203.333499999999  # GSM!M83 (-7574.428 -> -7203.333499999999)
$ws.Cells.Item(97, 8).Value = 340.3  # GSM!H97 (380.1111 -> 340.3)
$ws.Cells.Item(97, 9).Value = 211.44444  # GSM!I97 (255.66667 -> 211.44444)
$ws.Cells.Item(97, 10).Value = 1500  # GSM!J97 (629 -> 1500)
$ws.Cells.Item(97, 11).Value = 211.44444  # GSM!K97 (255.66667 -> 211.44444)
$ws.Cells.Item(97, 12).Value = 1500  # GSM!L97 (629 -> 1500)
$ws.Cells.Item(97, 13).Value = 284.55556  # GSM!M97 (240.33333 -> 284.55556)
$ws.Cells.Item(97, 14).Value = -2492  # GSM!N97 (-1621 -> -2492)
$ws.Cells.Item(102, 8).Value = 2360.1538  # GSM!H102 (3477.2856 -> 2360.1538)
$ws.Cells.Item(102, 10).Value = 3832.3333  # GSM!J102 (7374.25 -> 3832.3333)
$ws.Cells.Item(102, 12).Value = 3832.3333  # GSM!L102 (7374.25 -> 3832.3333)
$ws.Cells.Item(102, 14).Value = -7076.3333  # GSM!N102 (-10618.25 -> -7076.3333)
$ws.Cells.Item(122, 8).Value = 8040.25  # GSM!H122 (8044.375 -> 8040.25)
$ws.Cells.Item(122, 9).Value = 5665.0557  # GSM!I122 (5670.5557 -> 5665.0557)
$ws.Cells.Item(122, 11).Value = 16995.1671  # GSM!K122 (17011.6671 -> 16995.1671)
$ws.Cells.Item(122, 13).Value = -14545.1671  # GSM!M122 (-14561.6671 -> -14545.1671)
$ws.Cells.Item(141, 8).Value = 28707.2  # GSM!H141 (34574.6 -> 28707.2)
$ws.Cells.Item(141, 10).Value = 28707.2  # GSM!J141 (34574.6 -> 28707.2)
$ws.Cells.Item(141, 12).Value = 28707.2  # GSM!L141 (34574.6 -> 28707.2)
$ws.Cells.Item(141, 14).Value = -39067.2  # GSM!N141 (-44934.6 -> -39067.2)

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 33832.668  # LTW!H6 (43933.25 -> 33832.668)
$ws.Cells.Item(6, 10).Value = 47239.5  # LTW!J6 (56238 -> 47239.5)
$ws.Cells.Item(6, 12).Value = 47239.5  # LTW!L6 (56238 -> 47239.5)
$ws.Cells.Item(6, 14).Value = -47463.5  # LTW!N6 (-56462 -> -47463.5)
$ws.Cells.Item(22, 8).Value = 6246.75  # LTW!H22 (6248 -> 6246.75)
$ws.Cells.Item(22, 9).Value = 5845.6665  # LTW!I22 (5847.3335 -> 5845.6665)
$ws.Cells.Item(22, 11).Value = 5845.6665  # LTW!K22 (5847.3335 -> 5845.6665)
$ws.Cells.Item(22, 13).Value = -5550.6665  # LTW!M22 (-5552.3335 -> -5550.6665)
$ws.Cells.Item(27, 8).Value = 6246.75  # LTW!H27 (6248 -> 6246.75)
$ws.Cells.Item(27, 9).Value = 5845.6665  # LTW!I27 (5847.3335 -> 5845.6665)
$ws.Cells.Item(27, 11).Value = 5845.6665  # LTW!K27 (5847.3335 -> 5845.6665)
$ws.Cells.Item(27, 13).Value = -5738.6665  # LTW!M27 (-5740.3335 -> -5738.6665)
$ws.Cells.Item(40, 8).Value = 4713.4  # LTW!H40 (4275.5654 -> 4713.4)
$ws.Cells.Item(40, 9).Value = 3558.087  # LTW!I40 (3651.7273 -> 3558.087)
$ws.Cells.Item(40, 10).Value = 17999.5  # LTW!J40 (18000 -> 17999.5)
$ws.Cells.Item(40, 11).Value = 3558.087  # LTW!K40 (3651.7273 -> 3558.087)
$ws.Cells.Item(40, 12).Value = 17999.5  # LTW!L40 (18000 -> 17999.5)
$ws.Cells.Item(40, 13).Value = -3422.087  # LTW!M40 (-3515.7273 -> -3422.087)
$ws.Cells.Item(40, 14).Value = -18271.5  # LTW!N40 (-18272 -> -18271.5)
$ws.Cells.Item(46, 8).Value = 4322.1816  # LTW!H46 (4654.6 -> 4322.1816)
$ws.Cells.Item(46, 9).Value = 998.5  # LTW!I46 (999 -> 998.5)
$ws.Cells.Item(46, 11).Value = 998.5  # LTW!K46 (999 -> 998.5)
$ws.Cells.Item(46, 13).Value = -810.5  # LTW!M46 (-811 -> -810.5)
$ws.Cells.Item(93, 8).Value = 1232.0178  # LTW!H93 (1235.4108 -> 1232.0178)
$ws.Cells.Item(93, 9).Value = 987.7646999999999  # LTW!I93 (993.35297 -> 987.7646999999999)
$ws.Cells.Item(93, 11).Value = 987.7646999999999  # LTW!K93 (993.35297 -> 987.7646999999999)
$ws.Cells.Item(93, 13).Value = 260.2353000000001  # LTW!M93 (254.64703 -> 260.2353000000001)
$ws.Cells.Item(100, 8).Value = 5354.8335  # LTW!H100 (5506.1816 -> 5354.8335)
$ws.Cells.Item(100, 9).Value = 4425.5  # LTW!I100 (4507.222 -> 4425.5)
$ws.Cells.Item(100, 11).Value = 4425.5  # LTW!K100 (4507.222 -> 4425.5)
$ws.Cells.Item(100, 13).Value = -3884.5  # LTW!M100 (-3966.222 -> -3884.5)
$ws.Cells.Item(122, 8).Value = 6391.048  # LTW!H122 (6641.409 -> 6391.048)
$ws.Cells.Item(122, 10).Value = 12397.333  # LTW!J122 (12326.143 -> 12397.333)
$ws.Cells.Item(122, 12).Value = 37191.999  # LTW!L122 (36978.429 -> 37191.999)
$ws.Cells.Item(122, 14).Value = -42091.999  # LTW!N122 (-41878.429 -> -42091.999)
$ws.Cells.Item(136, 8).Value = 6846.8423  # LTW!H136 (7700.4375 -> 6846.8423)
$ws.Cells.Item(136, 9).Value = 3161.8293  # LTW!I136 (3587.8333 -> 3161.8293)
$ws.Cells.Item(136, 10).Value = 9644.723  # LTW!J136 (10168 -> 9644.723)
$ws.Cells.Item(136, 11).Value = 9485.4879  # LTW!K136 (10763.4999 -> 9485.4879)
$ws.Cells.Item(136, 12).Value = 28934.169  # LTW!L136 (30504 -> 28934.169)
$ws.Cells.Item(136, 13).Value = -6935.4879  # LTW!M136 (-8213.499899999999 -> -6935.4879)
$ws.Cells.Item(136, 14).Value = -34034.169  # LTW!N136 (-35604 -> -34034.169)

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 33834.78  # WVR!H81 (33836.156 -> 33834.78)
$ws.Cells.Item(81, 9).Value = 68345.8  # WVR!I81 (68348.734 -> 68345.8)
$ws.Cells.Item(81, 11).Value = 136691.6  # WVR!K81 (136697.468 -> 136691.6)
$ws.Cells.Item(81, 13).Value = -135630.6  # WVR!M81 (-135636.468 -> -135630.6)
$ws.Cells.Item(84, 8).Value = 33834.78  # WVR!H84 (33836.156 -> 33834.78)
$ws.Cells.Item(84, 9).Value = 68345.8  # WVR!I84 (68348.734 -> 68345.8)
$ws.Cells.Item(84, 11).Value = 683458  # WVR!K84 (683487.34 -> 683458)
$ws.Cells.Item(84, 13).Value = -678154  # WVR!M84 (-678183.34 -> -678154)
$ws.Cells.Item(107, 8).Value = 2928.8572  # WVR!H107 (2954.1428 -> 2928.8572)
$ws.Cells.Item(107, 9).Value = 3627  # WVR!I107 (3307.6 -> 3627)
$ws.Cells.Item(107, 10).Value = 1998  # WVR!J107 (2070.5 -> 1998)
$ws.Cells.Item(107, 11).Value = 10881  # WVR!K107 (9922.799999999999 -> 10881)
$ws.Cells.Item(107, 12).Value = 5994  # WVR!L107 (6211.5 -> 5994)
$ws.Cells.Item(107, 13).Value = -8961  # WVR!M107 (-8002.799999999999 -> -8961)
$ws.Cells.Item(107, 14).Value = -9834  # WVR!N107 (-10051.5 -> -9834)
$ws.Cells.Item(132, 8).Value = 3940.25  # WVR!H132 (3288.1538 -> 3940.25)
$ws.Cells.Item(132, 9).Value = 2378.0557  # WVR!I132 (2062.1667 -> 2378.0557)
$ws.Cells.Item(132, 11).Value = 7134.1671  # WVR!K132 (6186.500100000001 -> 7134.1671)
$ws.Cells.Item(132, 13).Value = -4604.1671  # WVR!M132 (-3656.500100000001 -> -4604.1671)
$ws.Cells.Item(140, 8).Value = 56752  # WVR!H140 (85485.2 -> 56752)
$ws.Cells.Item(140, 10).Value = 50940.75  # WVR!J140 (86857.25 -> 50940.75)
$ws.Cells.Item(140, 12).Value = 50940.75  # WVR!L140 (86857.25 -> 50940.75)
$ws.Cells.Item(140, 14).Value = -61300.75  # WVR!N140 (-97217.25 -> -61300.75)
